# Apply the "15-dec" column insertion to the "Prix Spot" sheet, and append
# two new date rows (2025-12-13, 2025-12-14) to the "Gaz" and "CO2" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Prix Spot" — insert a new column before ER (column 148) so
# everything from ER:FV shifts right to ES:FW, then fill the freed ER
# column with the new "15-dec" header and "-" placeholders for the data
# rows (2-25), matching the existing pattern used for empty days.
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Columns.Item(148).Insert()

$wsPrix.Range("ER1").Value2 = "15-dec"
$wsPrix.Range("ER2:ER25").Value2 = "-"

# ---------------------------------------------------------------------
# Sheets 2 & 3: "Gaz" and "CO2" — append two new rows for 2025-12-13 and
# 2025-12-14, continuing the existing list of dates in column A / prices
# in column B. The date cells must stay plain text (matching the rest of
# the column) rather than being auto-converted to a date serial, so we
# temporarily force a text number format, assign, then clear the format
# again so the resulting cell style matches the surrounding cells.
# ---------------------------------------------------------------------
function Add-DateRow {
    param($ws, $rowNum, $dateText, $value)

    $ws.Range("A$rowNum").NumberFormat = "@"
    $ws.Range("A$rowNum").Value2 = $dateText
    $ws.Range("A$rowNum").ClearFormats()
    $ws.Range("B$rowNum").Value2 = $value
}

$wsGaz = $wb.Worksheets.Item("Gaz")
Add-DateRow $wsGaz 178 "2025-12-13" 26.075
Add-DateRow $wsGaz 179 "2025-12-14" 26.075

$wsCO2 = $wb.Worksheets.Item("CO2")
Add-DateRow $wsCO2 178 "2025-12-13" 84.09999999999999
Add-DateRow $wsCO2 179 "2025-12-14" 84.09999999999999
